$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D remain plain text,
# matching the original workbook where these cells are stored as text.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.116.80'
$ws.Range('E2').Value = '  -1.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.257.22'
$ws.Range('E3').Value = '  -3.34%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '298.24'
$ws.Range('E5').Value = '  -2.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.35'
$ws.Range('E6').Value = '  -5.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.498'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  -3.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.87'
$ws.Range('E10').Value = '  -6.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0785'
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '48.23'
$ws.Range('E12').Value = '  -7.42%  '
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.63'
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.607.21'
$ws.Range('E15').Value = '  -3.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.33'
$ws.Range('E16').Value = '  -2.93%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.327.86'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.775'
$ws.Range('E18').Value = '  -2.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.093.13'
$ws.Range('E19').Value = '  -1.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0888'
$ws.Range('E20').Value = '  -2.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.40'
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.00'
$ws.Range('E22').Value = '  -3.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.38'
$ws.Range('E23').Value = '  -2.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '232.85'
$ws.Range('E24').Value = '  -1.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.91'
$ws.Range('E25').Value = '  -4.26%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.46'
$ws.Range('E26').Value = '  -3.88%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.78'
$ws.Range('E28').Value = '  -5.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.37'
$ws.Range('E29').Value = '  +4.38%  '
$ws.Range('E30').Value = '  -12.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.61'
$ws.Range('E31').Value = '  -3.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.03'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.93'
$ws.Range('E34').Value = '  -3.70%  '
$ws.Range('E35').Value = '  -4.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0692'
$ws.Range('E36').Value = '  -5.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.35'
$ws.Range('E37').Value = '  -5.16%  '
$ws.Range('E38').Value = '  -5.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.96'
$ws.Range('E39').Value = '  -8.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0983'
$ws.Range('E40').Value = '  -4.92%  '
$ws.Range('E41').Value = '  -3.14%  '
$ws.Range('E42').Value = '  -8.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.41'
$ws.Range('E43').Value = '  +2.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.935.62'
$ws.Range('E44').Value = '  -4.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0279'
$ws.Range('E45').Value = '  -2.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.29'
$ws.Range('E46').Value = '  -7.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.53'
$ws.Range('E47').Value = '  -7.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.77'
$ws.Range('E48').Value = '  -5.08%  '
$ws.Range('E49').Value = '  -3.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.482.82'
$ws.Range('E50').Value = '  -2.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.19'
